# Auto-generated: apply scheduled-runner price/profit refresh to Goblin_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2248.4707
$ws.Range("I33").Value = 372.22223
$ws.Range("J33").Value = 4359.25
$ws.Range("K33").Value = 372.22223
$ws.Range("L33").Value = 4359.25
$ws.Range("M33").Value = -143.22223
$ws.Range("N33").Value = -4817.25
$ws.Range("H62").Value = 5661.6665
$ws.Range("J62").Value = 8899.777
$ws.Range("L62").Value = 8899.777
$ws.Range("N62").Value = -10147.777
$ws.Range("H65").Value = 5661.6665
$ws.Range("J65").Value = 8899.777
$ws.Range("L65").Value = 44498.885
$ws.Range("N65").Value = -50738.885
$ws.Range("H86").Value = 4667
$ws.Range("I86").Value = 3003
$ws.Range("J86").Value = 5499
$ws.Range("K86").Value = 3003
$ws.Range("L86").Value = 5499
$ws.Range("M86").Value = -1880
$ws.Range("N86").Value = -7745
$ws.Range("H89").Value = 4667
$ws.Range("I89").Value = 3003
$ws.Range("J89").Value = 5499
$ws.Range("K89").Value = 15015
$ws.Range("L89").Value = 27495
$ws.Range("M89").Value = -9399
$ws.Range("N89").Value = -38727
$ws.Range("H98").Value = 1875
$ws.Range("I98").Value = 1815.909
$ws.Range("K98").Value = 1815.909
$ws.Range("M98").Value = -317.9090000000001
$ws.Range("H106").Value = 3474.2222
$ws.Range("I106").Value = 2879.8333
$ws.Range("J106").Value = 4663
$ws.Range("K106").Value = 2879.8333
$ws.Range("L106").Value = 4663
$ws.Range("M106").Value = -2248.8333
$ws.Range("N106").Value = -5925
$ws.Range("H122").Value = 1875
$ws.Range("I122").Value = 1815.909
$ws.Range("K122").Value = 5447.727000000001
$ws.Range("M122").Value = -2997.727000000001
$ws.Range("H138").Value = 3590961.5
$ws.Range("J138").Value = 6177478
$ws.Range("L138").Value = 18532434
$ws.Range("N138").Value = -18542714
$ws.Range("H141").Value = 4887.5884
$ws.Range("I141").Value = 4887.5884
$ws.Range("K141").Value = 14662.7652
$ws.Range("M141").Value = -9482.765199999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4916.6567
$ws.Range("I32").Value = 3022.9836
$ws.Range("K32").Value = 3022.9836
$ws.Range("M32").Value = -2735.9836
$ws.Range("H92").Value = 46250
$ws.Range("J92").Value = 46250
$ws.Range("L92").Value = 46250
$ws.Range("N92").Value = -51242
$ws.Range("H102").Value = 3963.2903
$ws.Range("I102").Value = 1863.5652
$ws.Range("K102").Value = 1863.5652
$ws.Range("M102").Value = -241.5652
$ws.Range("H132").Value = 6120
$ws.Range("I132").Value = 6064.2085
$ws.Range("J132").Value = 6454.75
$ws.Range("K132").Value = 18192.6255
$ws.Range("L132").Value = 19364.25
$ws.Range("M132").Value = -15662.6255
$ws.Range("N132").Value = -24424.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 38466450
$ws.Range("I86").Value = 6265
$ws.Range("K86").Value = 6265
$ws.Range("M86").Value = -5142
$ws.Range("H89").Value = 38466450
$ws.Range("I89").Value = 6265
$ws.Range("K89").Value = 31325
$ws.Range("M89").Value = -25709
$ws.Range("H107").Value = 3789.476
$ws.Range("I107").Value = 776.7
$ws.Range("J107").Value = 6528.364
$ws.Range("K107").Value = 776.7
$ws.Range("L107").Value = 6528.364
$ws.Range("M107").Value = 1143.3
$ws.Range("N107").Value = -10368.364
$ws.Range("H134").Value = 964226.2
$ws.Range("I134").Value = 2625.3157
$ws.Range("K134").Value = 7875.9471
$ws.Range("M134").Value = -5340.9471
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1222.6923
$ws.Range("I22").Value = 372.7143
$ws.Range("K22").Value = 372.7143
$ws.Range("M22").Value = -22.71429999999998
$ws.Range("H105").Value = 4554.6665
$ws.Range("I105").Value = 3582
$ws.Range("K105").Value = 3582
$ws.Range("M105").Value = -1835
$ws.Range("H107").Value = 111952.11
$ws.Range("I107").Value = 251212.25
$ws.Range("J107").Value = 544
$ws.Range("K107").Value = 251212.25
$ws.Range("L107").Value = 544
$ws.Range("M107").Value = -249292.25
$ws.Range("N107").Value = -4384
$ws.Range("H132").Value = 2278.111
$ws.Range("I132").Value = 2534
$ws.Range("J132").Value = 1766.3334
$ws.Range("K132").Value = 7602
$ws.Range("L132").Value = 5299.0002
$ws.Range("M132").Value = -5072
$ws.Range("N132").Value = -10359.0002
$ws.Range("H134").Value = 2521.389
$ws.Range("I134").Value = 1899.0625
$ws.Range("K134").Value = 5697.1875
$ws.Range("M134").Value = -3162.1875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 86.85
$ws.Range("I4").Value = 87.121216
$ws.Range("K4").Value = 261.363648
$ws.Range("M4").Value = -149.363648
$ws.Range("H7").Value = 236.7
$ws.Range("I7").Value = 171.83333
$ws.Range("K7").Value = 515.49999
$ws.Range("M7").Value = -403.49999
$ws.Range("H12").Value = 356.83334
$ws.Range("I12").Value = 229.66667
$ws.Range("J12").Value = 484
$ws.Range("K12").Value = 689.00001
$ws.Range("L12").Value = 1452
$ws.Range("M12").Value = -516.00001
$ws.Range("N12").Value = -1798
$ws.Range("H50").Value = 20050.154
$ws.Range("I50").Value = 357.83334
$ws.Range("J50").Value = 36929.285
$ws.Range("K50").Value = 1073.50002
$ws.Range("L50").Value = 110787.855
$ws.Range("M50").Value = -592.50002
$ws.Range("N50").Value = -111749.855
$ws.Range("H53").Value = 20050.154
$ws.Range("I53").Value = 357.83334
$ws.Range("J53").Value = 36929.285
$ws.Range("K53").Value = 1073.50002
$ws.Range("L53").Value = 110787.855
$ws.Range("M53").Value = -592.50002
$ws.Range("N53").Value = -111749.855
$ws.Range("H121").Value = 52160.95
$ws.Range("J121").Value = 3519.3635
$ws.Range("L121").Value = 10558.0905
$ws.Range("N121").Value = -13178.0905
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 750012500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 750012500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 750012500
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -750013872
$ws.Range("H65").Value = 750012500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 750012500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 2250037500
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -2250044364
$ws.Range("H102").Value = 1712
$ws.Range("I102").Value = 1461.3334
$ws.Range("K102").Value = 1461.3334
$ws.Range("M102").Value = 160.6666
$ws.Range("H122").Value = 1907.1852
$ws.Range("I122").Value = 1537.4762
$ws.Range("K122").Value = 4612.4286
$ws.Range("M122").Value = -2162.4286
$ws.Range("H132").Value = 43480988
$ws.Range("I132").Value = 58824536
$ws.Range("K132").Value = 176473608
$ws.Range("M132").Value = -176471078
$ws.Range("H136").Value = 19489.143
$ws.Range("J136").Value = 19489.143
$ws.Range("L136").Value = 58467.429
$ws.Range("N136").Value = -63567.429
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1738.1538
$ws.Range("J46").Value = 2049.6
$ws.Range("L46").Value = 2049.6
$ws.Range("N46").Value = -2425.6
$ws.Range("H62").Value = 400096000
$ws.Range("J62").Value = 400096000
$ws.Range("L62").Value = 400096000
$ws.Range("N62").Value = -400097248
$ws.Range("H65").Value = 400096000
$ws.Range("J65").Value = 400096000
$ws.Range("L65").Value = 1200288000
$ws.Range("N65").Value = -1200294240
$ws.Range("H95").Value = 39835.75
$ws.Range("J95").Value = 39835.75
$ws.Range("L95").Value = 39835.75
$ws.Range("N95").Value = -45327.75
$ws.Range("H96").Value = 22000
$ws.Range("J96").Value = 22000
$ws.Range("L96").Value = 22000
$ws.Range("N96").Value = -27492
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 35271
$ws.Range("J69").Value = 35271
$ws.Range("L69").Value = 35271
$ws.Range("N69").Value = -36769
$ws.Range("H72").Value = 35271
$ws.Range("J72").Value = 35271
$ws.Range("L72").Value = 105813
$ws.Range("N72").Value = -113301
$ws.Range("H126").Value = 9261828
$ws.Range("I126").Value = 2483.2222
$ws.Range("K126").Value = 7449.6666
$ws.Range("M126").Value = -4979.6666
$ws.Range("H132").Value = 8779048
$ws.Range("I132").Value = 9266730
$ws.Range("J132").Value = 764
$ws.Range("K132").Value = 27800190
$ws.Range("L132").Value = 2292
$ws.Range("M132").Value = -27797660
$ws.Range("N132").Value = -7352
$ws.Range("H136").Value = 7101.364
$ws.Range("I136").Value = 5981
$ws.Range("J136").Value = 8445.8
$ws.Range("K136").Value = 17943
$ws.Range("L136").Value = 25337.4
$ws.Range("M136").Value = -15393
$ws.Range("N136").Value = -30437.4
